$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'249.81"
$ws.Range("G2").Value = "'11"

# Row 3
$ws.Range("D3").Value = "'23.84"
$ws.Range("G3").Value = "'11"

# Row 4
$ws.Range("D4").Value = "'6.031"
$ws.Range("G4").Value = "'11"

# Row 5
$ws.Range("D5").Value = "'0.05978"
$ws.Range("G5").Value = "'11"

# Row 6
$ws.Range("D6").Value = "'3.424"
$ws.Range("G6").Value = "'11"

# Row 7
$ws.Range("D7").Value = "'6.559"
$ws.Range("G7").Value = "'11"

# Row 8
$ws.Range("G8").Value = "'11"

# Row 9
$ws.Range("D9").Value = "'0.7973"
$ws.Range("G9").Value = "'11"

# Row 10
$ws.Range("B10").Value = "'One"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.0006106"
$ws.Range("E10").Value = "'9OneONE"
$ws.Range("G10").Value = "'11"

# Row 11
$ws.Range("B11").Value = "'WazirX"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1507"
$ws.Range("E11").Value = "'10WazirXWRX"
$ws.Range("G11").Value = "'11"

# Row 12
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07894"
$ws.Range("E12").Value = "'11MandalaExchangeTokenMDX"
$ws.Range("G12").Value = "'11"

# Row 13
$ws.Range("B13").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03347"
$ws.Range("E13").Value = "'12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").Value = "'11"

# Row 14
$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03060"
$ws.Range("E14").Value = "'13BitrueCoinBTR"
$ws.Range("G14").Value = "'11"

# Row 15
$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09279"
$ws.Range("E15").Value = "'14BitMartTokenBMX"
$ws.Range("G15").Value = "'11"

# Row 16
$ws.Range("B16").Value = "'MCDex"
$ws.Range("C16").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.565"
$ws.Range("E16").Value = "'15MCDexMCB"
$ws.Range("G16").Value = "'11"

# Row 17
$ws.Range("B17").Value = "'BitForexToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001668"
$ws.Range("E17").Value = "'16BitForexTokenBF"
$ws.Range("G17").Value = "'11"

# Row 18
$ws.Range("B18").Value = "'CoinExToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04776"
$ws.Range("E18").Value = "'17CoinExTokenCET"
$ws.Range("G18").Value = "'11"

# Row 19
$ws.Range("D19").Value = "'0.006230"
$ws.Range("G19").Value = "'11"

# Row 20
$ws.Range("D20").Value = "'0.005682"
$ws.Range("G20").Value = "'11"

# Row 21
$ws.Range("D21").Value = "'0.001073"
$ws.Range("G21").Value = "'11"

# Row 22
$ws.Range("G22").Value = "'11"

# Row 23
$ws.Range("G23").Value = "'11"

# Row 24
$ws.Range("D24").Value = "'2.206"
$ws.Range("G24").Value = "'11"

# Row 25
$ws.Range("D25").Value = "'0.3350"
$ws.Range("G25").Value = "'11"

# Row 26
$ws.Range("G26").Value = "'11"

# Row 27
$ws.Range("D27").Value = "'0.0006492"
$ws.Range("G27").Value = "'11"

# Row 28
$ws.Range("G28").Value = "'11"

# Row 29
$ws.Range("G29").Value = "'11"

# Row 30
$ws.Range("G30").Value = "'11"

# Row 31
$ws.Range("G31").Value = "'11"

# Row 32
$ws.Range("G32").Value = "'11"

# Row 33
$ws.Range("G33").Value = "'11"

# Row 34
$ws.Range("G34").Value = "'11"

# Row 35
$ws.Range("G35").Value = "'11"

# Row 36
$ws.Range("G36").Value = "'11"

# Row 37
$ws.Range("G37").Value = "'11"

# Row 38
$ws.Range("G38").Value = "'11"

# Row 39
$ws.Range("G39").Value = "'11"

# Row 40
$ws.Range("D40").Value = "'0.04431"
$ws.Range("G40").Value = "'11"

# Row 41
$ws.Range("D41").Value = "'0.007043"
$ws.Range("G41").Value = "'11"

# Row 42
$ws.Range("B42").Value = "'CEJI"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003610"
$ws.Range("E42").Value = "'41CEJICEJI"
$ws.Range("G42").Value = "'11"

# Row 43
$ws.Range("B43").Value = "'BKEXToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1066"
$ws.Range("E43").Value = "'42BKEXTokenBKK"
$ws.Range("G43").Value = "'11"

# Row 44
$ws.Range("D44").Value = "'0.009228"
$ws.Range("G44").Value = "'11"

# Row 45
$ws.Range("G45").Value = "'11"

# Row 46
$ws.Range("D46").Value = "'0.00005890"
$ws.Range("G46").Value = "'11"

# Row 47
$ws.Range("G47").Value = "'11"

# Row 48
$ws.Range("D48").Value = "'0.9932"
$ws.Range("G48").Value = "'11"

# Row 49
$ws.Range("D49").Value = "'0.09649"
$ws.Range("E49").Value = "'48BOLOBOLOWorstin24h"
$ws.Range("G49").Value = "'11"

# Row 50
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("G50").Value = "'11"

# Row 51
$ws.Range("D51").Value = "'0.01013"
$ws.Range("G51").Value = "'11"
